$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price values remain stored as text (matching source format)
$textFormatCells = @('D5', 'D6', 'D20', 'D21', 'D26', 'D27', 'D31', 'D32', 'D35', 'D37', 'D40', 'D41', 'D46', 'D47', 'D51')
foreach ($cellRef in $textFormatCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply updated cell values from the crypto data refresh
$ws.Range('D2').Value = '70.611.30'
$ws.Range('E2').Value = '  -1.63%  '
$ws.Range('D3').Value = '2.536.13'
$ws.Range('E3').Value = '  -5.06%  '
$ws.Range('D5').Value = '576.82'
$ws.Range('E5').Value = '  -3.44%  '
$ws.Range('D6').Value = '170.10'
$ws.Range('E6').Value = '  -2.79%  '
$ws.Range('E8').Value = '  -2.58%  '
$ws.Range('E9').Value = '  -0.79%  '
$ws.Range('D10').Value = '2.536.21'
$ws.Range('E10').Value = '  -5.01%  '
$ws.Range('E11').Value = '  -0.21%  '
$ws.Range('E12').Value = '  -3.27%  '
$ws.Range('E13').Value = '  -3.49%  '
$ws.Range('D14').Value = '2.998.03'
$ws.Range('E14').Value = '  -5.44%  '
$ws.Range('D15').Value = '70.458.08'
$ws.Range('E15').Value = '  -1.72%  '
$ws.Range('E16').Value = '  -2.34%  '
$ws.Range('E17').Value = '  -4.27%  '
$ws.Range('D18').Value = '2.536.05'
$ws.Range('E18').Value = '  -4.99%  '
$ws.Range('E19').Value = '  -4.70%  '
$ws.Range('D20').Value = '360.32'
$ws.Range('D21').Value = '7.37'
$ws.Range('E21').Value = '  -10.02%  '
$ws.Range('E22').Value = '  -5.52%  '
$ws.Range('E23').Value = '  -1.82%  '
$ws.Range('E25').Value = '  -2.96%  '
$ws.Range('D26').Value = '4.09'
$ws.Range('E26').Value = '  -5.45%  '
$ws.Range('D27').Value = '9.26'
$ws.Range('E27').Value = '  -4.98%  '
$ws.Range('D28').Value = '2.666.34'
$ws.Range('E28').Value = '  -5.34%  '
$ws.Range('E29').Value = '  +0.10%  '
$ws.Range('D30').Value = '0.0₃0921'
$ws.Range('E30').Value = '  -4.58%  '
$ws.Range('D31').Value = '7.88'
$ws.Range('E31').Value = '  -2.08%  '
$ws.Range('D32').Value = '483.69'
$ws.Range('E32').Value = '  -3.57%  '
$ws.Range('E33').Value = '  -1.45%  '
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('E36').Value = '  +6.03%  '
$ws.Range('D37').Value = '157.38'
$ws.Range('E37').Value = '  -3.69%  '
$ws.Range('E39').Value = '  -1.28%  '
$ws.Range('B40').Value = 'ImmutableX'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D40').Value = '1.31'
$ws.Range('E40').Value = '  -4.52%  '
$ws.Range('B41').Value = 'USDe'
$ws.Range('C41').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('E42').Value = '  -5.77%  '
$ws.Range('E43').Value = '  -4.90%  '
$ws.Range('E44').Value = '  -4.58%  '
$ws.Range('E45').Value = '  -4.09%  '
$ws.Range('D46').Value = '38.38'
$ws.Range('E46').Value = '  -2.76%  '
$ws.Range('D47').Value = '144.41'
$ws.Range('E47').Value = '  -7.22%  '
$ws.Range('E48').Value = '  -4.84%  '
$ws.Range('E49').Value = '  -5.86%  '
$ws.Range('E50').Value = '  -6.50%  '
$ws.Range('D51').Value = '0.595'
$ws.Range('E51').Value = '  -1.38%  '
